# Add 5 new vocabulary rows to the ENGLISH sheet and 5 new quotes to the
# NOTES sheet (new words/quotes captured during the 2021-11-14 16:xx study
# session).

$wb = $excel.ActiveWorkbook
$wsEnglish = $wb.Worksheets.Item("ENGLISH")
$wsNotes = $wb.Worksheets.Item("NOTES")

# Helper: write a genuinely empty-string value into a cell (matches the
# source export, which stores blank Definition/Antonyms as empty strings
# rather than leaving the cell absent). A lone quote-prefix character
# evaluates to an empty text value; ClearFormats() strips the quote-prefix
# cell style back off again so no visible formatting change remains.
function Set-EmptyString($cell) {
    $cell.Value = "'"
    $cell.ClearFormats()
}

# New vocabulary entries (Word, Definition, Synonyms, Antonyms, Correct
# answer count, Created at)
$newWords = @(
    @("erratic", "", "unpredictable", "", 0, "2021-11-14 16:37:51.305571"),
    @("entail", "involve (smth) as a necessary or inevitable part or consequence", "necessitate", "", 0, "2021-11-14 16:39:12.99207"),
    @("disengage", "", "remove;withdraw", "", 0, "2021-11-14 16:45:31.492612"),
    @("recuperation", "recovery from illness or exertion", "recovery", "", 0, "2021-11-14 16:46:27.758408"),
    @("nurture", "care for and protect (someone or smth) while they are growing", "cultivate", "", 0, "2021-11-14 16:48:04.864399")
)

$startRow = 76
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = $startRow + $i
    $data = $newWords[$i]

    $wsEnglish.Cells.Item($row, 1).Value = $data[0]

    if ($data[1] -eq "") {
        Set-EmptyString($wsEnglish.Cells.Item($row, 2))
    } else {
        $wsEnglish.Cells.Item($row, 2).Value = $data[1]
    }

    $wsEnglish.Cells.Item($row, 3).Value = $data[2]

    if ($data[3] -eq "") {
        Set-EmptyString($wsEnglish.Cells.Item($row, 4))
    } else {
        $wsEnglish.Cells.Item($row, 4).Value = $data[3]
    }

    $wsEnglish.Cells.Item($row, 5).Value = $data[4]
    $wsEnglish.Cells.Item($row, 6).Value = $data[5]
}

# New quotes appended to the NOTES sheet
$newNotes = @(
    "Each night, when I go to sleep, I die. And the next morning, when I wake up, I am reborn",
    "If the answer isn't a definite yes then it should be a no",
    "If it isn't a clear yes, then it's a clear no",
    "To follow, without halt, one aim: there is the secret to success",
    "Half of the troubles of this life can be traced to saying yes too quickly and not saying no soon enough"
)

$notesStartRow = 21
for ($i = 0; $i -lt $newNotes.Length; $i++) {
    $row = $notesStartRow + $i
    $wsNotes.Cells.Item($row, 1).Value = $newNotes[$i]
}
